$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update 想去人数 (column F) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 66
$ws1.Range("F8").Value = 3874
$ws1.Range("F10").Value = 4563
$ws1.Range("F12").Value = 1150

# Sheet "全部类型" (sheet4): update 想去人数 (column F) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 66
$ws4.Range("F9").Value = 3874
$ws4.Range("F11").Value = 4563
$ws4.Range("F13").Value = 1150
